# Update the "想去人数" (F column) figures for the 展览 (Exhibitions) sheet
# and the identical 全部类型 (All Types) sheet, as scraped at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Row -> new value for column F on both affected sheets.
$updates = @{
    4  = 1400
    5  = 473
    6  = 210
    9  = 132
    11 = 343
    13 = 1829
    15 = 113
    16 = 185
    17 = 717
    20 = 4368
    22 = 317
    23 = 1180
    26 = 735
    28 = 383
    30 = 193
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
